$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.109.63'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '2.929.97'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.33'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.10'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('E9').Value = '  +4.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  -0.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.442'
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '33.71'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').Value = '3.413.48'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '61.106.96'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.73'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = '2.929.21'
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '433.69'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.47'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '81.56'
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('E24').Value = '  +2.73%  '
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.86'
$ws.Range('E26').Value = '  -1.07%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  -1.42%  '
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.98'
$ws.Range('E30').Value = '  -0.96%  '
$ws.Range('E31').Value = '  +2.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.73'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '0.0₃0869'
$ws.Range('E34').Value = '  +1.79%  '
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.65'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.99'
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.60'
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '42.07'
$ws.Range('E41').Value = '  +5.05%  '
$ws.Range('E42').Value = '  -2.55%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0347'
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '372.61'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').Value = '2.710.12'
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '133.15'
$ws.Range('E46').Value = '  +1.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.88'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('E49').Value = '  -1.15%  '
$ws.Range('E50').Value = '  -1.27%  '
$ws.Range('E51').Value = '  -0.42%  '
